# Generate Report for Handoff
# The "d5ef8c1c-1440-43a9-9dca-75493d500fed" file has finished translation
# and is now ready for handoff. Update its status + handoff timestamps on
# the Overview sheet and on each per-language (zh-cn / de-de) detail sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the d5ef8c1c.md file -------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"          # zh-cn status
$overview.Range("C3").Value = "Ready for handoff"          # de-de status
$overview.Range("D3").Value = "2016-45-18 03:45:18"        # Latest Handoff Date

# --- zh-cn detail sheet: row 3 is the d5ef8c1c file -------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"               # Status
$zhcn.Range("E3").Value = "2016-03-18 03:45:16"             # Latest Handoff Datetime

# --- de-de detail sheet: row 3 is the d5ef8c1c file -------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"                # Status
$dede.Range("E3").Value = "2016-03-18 03:45:18"              # Latest Handoff Datetime
